$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1077.1428
$ws.Range("J17").Value = 1077.1428
$ws.Range("L17").Value = 3231.4284
$ws.Range("N17").Value = -3567.4284

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 61718.883
$ws.Range("J64").Value = 3180.111
$ws.Range("L64").Value = 3180.111
$ws.Range("N64").Value = -3676.111

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 61718.883
$ws.Range("J67").Value = 3180.111
$ws.Range("L67").Value = 3180.111
$ws.Range("N67").Value = -4896.111

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3072.25
$ws.Range("J74").Value = 3121.25
$ws.Range("L74").Value = 3121.25
$ws.Range("N74").Value = -4993.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3072.25
$ws.Range("J77").Value = 3121.25
$ws.Range("L77").Value = 15606.25
$ws.Range("N77").Value = -24966.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1847.1666
$ws.Range("I127").Value = 487.44446
$ws.Range("J127").Value = 2300.4075
$ws.Range("K127").Value = 1462.33338
$ws.Range("L127").Value = 6901.2225
$ws.Range("M127").Value = 3497.66662
$ws.Range("N127").Value = -16821.2225

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1496.2858
$ws.Range("I135").Value = 383.73334
$ws.Range("J135").Value = 2330.7
$ws.Range("K135").Value = 3453.60006
$ws.Range("L135").Value = 20976.3
$ws.Range("M135").Value = -918.6000599999998
$ws.Range("N135").Value = -26046.3

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1623.3934
$ws.Range("I138").Value = 1594.6666
$ws.Range("J138").Value = 1635.4186
$ws.Range("K138").Value = 4783.9998
$ws.Range("L138").Value = 4906.2558
$ws.Range("M138").Value = 356.0002000000004
$ws.Range("N138").Value = -15186.2558

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 60330.8
$ws.Range("J139").Value = 60330.8
$ws.Range("L139").Value = 60330.8
$ws.Range("N139").Value = -70610.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22025.467
$ws.Range("I32").Value = 5665.047
$ws.Range("J32").Value = 220687.72
$ws.Range("K32").Value = 5665.047
$ws.Range("L32").Value = 220687.72
$ws.Range("M32").Value = -5378.047
$ws.Range("N32").Value = -221261.72

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2163.9656
$ws.Range("I61").Value = 1709.1177
$ws.Range("J61").Value = 2808.3333
$ws.Range("K61").Value = 1709.1177
$ws.Range("L61").Value = 2808.3333
$ws.Range("M61").Value = -1497.1177
$ws.Range("N61").Value = -3232.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 64454.5
$ws.Range("I97").Value = 144254.28
$ws.Range("J97").Value = 2388
$ws.Range("K97").Value = 144254.28
$ws.Range("L97").Value = 2388
$ws.Range("M97").Value = -143758.28
$ws.Range("N97").Value = -3380

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2163.9656
$ws.Range("I136").Value = 1709.1177
$ws.Range("J136").Value = 2808.3333
$ws.Range("K136").Value = 5127.3531
$ws.Range("L136").Value = 8424.999899999999
$ws.Range("M136").Value = -2577.3531
$ws.Range("N136").Value = -13524.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 7500
$ws.Range("J9").Value = 7500
$ws.Range("L9").Value = 7500
$ws.Range("N9").Value = -7836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1035.6
$ws.Range("I94").Value = 1094.5
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 1094.5
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = -643.5
$ws.Range("N94").Value = -1702

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 85000
$ws.Range("J138").Value = 85000
$ws.Range("L138").Value = 85000
$ws.Range("N138").Value = -95280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3003375
$ws.Range("J4").Value = 671166.7
$ws.Range("L4").Value = 671166.7
$ws.Range("N4").Value = -671390.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 62720
$ws.Range("J100").Value = 62720
$ws.Range("L100").Value = 62720
$ws.Range("N100").Value = -64884

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4632.778
$ws.Range("I132").Value = 5137.846
$ws.Range("J132").Value = 3319.6
$ws.Range("K132").Value = 15413.538
$ws.Range("L132").Value = 9958.799999999999
$ws.Range("M132").Value = -12883.538
$ws.Range("N132").Value = -15018.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2490
$ws.Range("I134").Value = 980
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 2940
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -405
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 54749.5
$ws.Range("J140").Value = 54749.5
$ws.Range("L140").Value = 54749.5
$ws.Range("N140").Value = -65109.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1237.5
$ws.Range("J34").Value = 1574.75
$ws.Range("L34").Value = 4724.25
$ws.Range("N34").Value = -4892.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 19130
$ws.Range("I55").Value = 25575
$ws.Range("J55").Value = 14833.333
$ws.Range("K55").Value = 76725
$ws.Range("L55").Value = 44499.999
$ws.Range("M55").Value = -76548
$ws.Range("N55").Value = -44853.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1002.5
$ws.Range("I59").Value = 1005
$ws.Range("K59").Value = 3015
$ws.Range("M59").Value = -2475

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1252500
$ws.Range("J5").Value = 4500
$ws.Range("L5").Value = 4500
$ws.Range("N5").Value = -4724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 40001880
$ws.Range("I97").Value = 52633356
$ws.Range("J97").Value = 2195.3333
$ws.Range("K97").Value = 52633356
$ws.Range("L97").Value = 2195.3333
$ws.Range("M97").Value = -52632860
$ws.Range("N97").Value = -3187.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 219047.62
$ws.Range("J2").Value = 8333.333000000001
$ws.Range("L2").Value = 8333.333000000001
$ws.Range("N2").Value = -8557.333000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5436.4375
$ws.Range("I7").Value = 3758.2856
$ws.Range("K7").Value = 3758.2856
$ws.Range("M7").Value = -3646.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5436.4375
$ws.Range("I126").Value = 3758.2856
$ws.Range("K126").Value = 11274.8568
$ws.Range("M126").Value = -8804.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 9866.666999999999
$ws.Range("I4").Value = 10000
$ws.Range("J4").Value = 9600
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 9600
$ws.Range("M4").Value = -9887
$ws.Range("N4").Value = -9826

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7266.25
$ws.Range("J41").Value = 7266.25
$ws.Range("L41").Value = 7266.25
$ws.Range("N41").Value = -8046.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 55556480
$ws.Range("I96").Value = 90910104
$ws.Range("J96").Value = 783.8570999999999
$ws.Range("K96").Value = 90910104
$ws.Range("L96").Value = 783.8570999999999
$ws.Range("M96").Value = -90908731
$ws.Range("N96").Value = -3529.8571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 944.4
$ws.Range("I113").Value = 845
$ws.Range("J113").Value = 1129
$ws.Range("K113").Value = 2535
$ws.Range("L113").Value = 3387
$ws.Range("M113").Value = -365
$ws.Range("N113").Value = -7727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2247.077
$ws.Range("I126").Value = 2109.7273
$ws.Range("K126").Value = 6329.1819
$ws.Range("M126").Value = -3859.1819
